$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 7).Value = 28.479168
$ws.Cells.Item(2, 8).Value = 85.437504
$ws.Cells.Item(2, 9).Value = 0.4446244458164738
$ws.Cells.Item(2, 10).Value = 0.4446244458164738
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 13).Value = 87.94127933333334
$ws.Cells.Item(2, 14).Value = 263.823838
$ws.Cells.Item(2, 15).Value = 0.4109331243514438
$ws.Cells.Item(2, 16).Value = 0.4109331243514437
$ws.Cells.Item(2, 17).Value = 2504.494468268928
$ws.Cells.Item(2, 18).Value = 22540.45021442035
$ws.Cells.Item(2, 19).Value = 0.1827109126823928
$ws.Cells.Item(2, 20).Value = 0.1827109126823928

# Row 3
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 7).Value = 28.479168
$ws.Cells.Item(3, 8).Value = 85.437504
$ws.Cells.Item(3, 9).Value = 0.4446244458164738
$ws.Cells.Item(3, 10).Value = 0.4446244458164738
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 13).Value = 52.441971
$ws.Cells.Item(3, 14).Value = 157.325913
$ws.Cells.Item(3, 15).Value = 0.2450515065683088
$ws.Cells.Item(3, 16).Value = 0.2450515065683087
$ws.Cells.Item(3, 17).Value = 1493.503702360128
$ws.Cells.Item(3, 18).Value = 13441.53332124115
$ws.Cells.Item(3, 19).Value = 0.1089558903044263
$ws.Cells.Item(3, 20).Value = 0.1089558903044263

# Row 4
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 7).Value = 28.479168
$ws.Cells.Item(4, 8).Value = 85.437504
$ws.Cells.Item(4, 9).Value = 0.4446244458164738
$ws.Cells.Item(4, 10).Value = 0.4446244458164738
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 13).Value = 54.667459
$ws.Cells.Item(4, 14).Value = 164.002377
$ws.Cells.Item(4, 15).Value = 0.255450795093328
$ws.Cells.Item(4, 16).Value = 0.255450795093328
$ws.Cells.Item(4, 17).Value = 1556.883748994112
$ws.Cells.Item(4, 18).Value = 14011.95374094701
$ws.Cells.Item(4, 19).Value = 0.1135796682017486
$ws.Cells.Item(4, 20).Value = 0.1135796682017486

# Row 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 7).Value = 28.479168
$ws.Cells.Item(5, 8).Value = 85.437504
$ws.Cells.Item(5, 9).Value = 0.4446244458164738
$ws.Cells.Item(5, 10).Value = 0.4446244458164738
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 13).Value = 18.95316166666667
$ws.Cells.Item(5, 14).Value = 56.85948500000001
$ws.Cells.Item(5, 15).Value = 0.08856457398691947
$ws.Cells.Item(5, 16).Value = 0.08856457398691944
$ws.Cells.Item(5, 17).Value = 539.7702752361602
$ws.Cells.Item(5, 18).Value = 4857.932477125441
$ws.Cells.Item(5, 19).Value = 0.03937797462790616
$ws.Cells.Item(5, 20).Value = 0.03937797462790615

# Row 6
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 7).Value = 18.12667766666667
$ws.Cells.Item(6, 8).Value = 54.380033
$ws.Cells.Item(6, 9).Value = 0.2829985767855128
$ws.Cells.Item(6, 10).Value = 0.2829985767855128
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 13).Value = 87.94127933333334
$ws.Cells.Item(6, 14).Value = 263.823838
$ws.Cells.Item(6, 15).Value = 0.4109331243514438
$ws.Cells.Item(6, 16).Value = 0.4109331243514437
$ws.Cells.Item(6, 17).Value = 1594.083224069629
$ws.Cells.Item(6, 18).Value = 14346.74901662666
$ws.Cells.Item(6, 19).Value = 0.1162934893454828
$ws.Cells.Item(6, 20).Value = 0.1162934893454827

# Row 7
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 7).Value = 18.12667766666667
$ws.Cells.Item(7, 8).Value = 54.380033
$ws.Cells.Item(7, 9).Value = 0.2829985767855128
$ws.Cells.Item(7, 10).Value = 0.2829985767855128
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 13).Value = 52.441971
$ws.Cells.Item(7, 14).Value = 157.325913
$ws.Cells.Item(7, 15).Value = 0.2450515065683088
$ws.Cells.Item(7, 16).Value = 0.2450515065683087
$ws.Cells.Item(7, 17).Value = 950.5987045216812
$ws.Cells.Item(7, 18).Value = 8555.38834069513
$ws.Cells.Item(7, 19).Value = 0.06934922759797713
$ws.Cells.Item(7, 20).Value = 0.06934922759797713

# Row 8
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 7).Value = 18.12667766666667
$ws.Cells.Item(8, 8).Value = 54.380033
$ws.Cells.Item(8, 9).Value = 0.2829985767855128
$ws.Cells.Item(8, 10).Value = 0.2829985767855128
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 13).Value = 54.667459
$ws.Cells.Item(8, 14).Value = 164.002377
$ws.Cells.Item(8, 15).Value = 0.255450795093328
$ws.Cells.Item(8, 16).Value = 0.255450795093328
$ws.Cells.Item(8, 17).Value = 990.9394081487159
$ws.Cells.Item(8, 18).Value = 8918.454673338441
$ws.Cells.Item(8, 19).Value = 0.07229221145013949
$ws.Cells.Item(8, 20).Value = 0.07229221145013948

# Row 9
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 7).Value = 18.12667766666667
$ws.Cells.Item(9, 8).Value = 54.380033
$ws.Cells.Item(9, 9).Value = 0.2829985767855128
$ws.Cells.Item(9, 10).Value = 0.2829985767855128
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 13).Value = 18.95316166666667
$ws.Cells.Item(9, 14).Value = 56.85948500000001
$ws.Cells.Item(9, 15).Value = 0.08856457398691947
$ws.Cells.Item(9, 16).Value = 0.08856457398691944
$ws.Cells.Item(9, 17).Value = 343.5578522958896
$ws.Cells.Item(9, 18).Value = 3092.020670663006
$ws.Cells.Item(9, 19).Value = 0.02506364839191346
$ws.Cells.Item(9, 20).Value = 0.02506364839191345

# Row 10
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 7).Value = 11.513346
$ws.Cells.Item(10, 8).Value = 34.540038
$ws.Cells.Item(10, 9).Value = 0.179749460544048
$ws.Cells.Item(10, 10).Value = 0.179749460544048
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 13).Value = 87.94127933333334
$ws.Cells.Item(10, 14).Value = 263.823838
$ws.Cells.Item(10, 15).Value = 0.4109331243514438
$ws.Cells.Item(10, 16).Value = 0.4109331243514437
$ws.Cells.Item(10, 17).Value = 1012.498376647316
$ws.Cells.Item(10, 18).Value = 9112.485389825846
$ws.Cells.Item(10, 19).Value = 0.07386500742185223
$ws.Cells.Item(10, 20).Value = 0.07386500742185222

# Row 11
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 7).Value = 11.513346
$ws.Cells.Item(11, 8).Value = 34.540038
$ws.Cells.Item(11, 9).Value = 0.179749460544048
$ws.Cells.Item(11, 10).Value = 0.179749460544048
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 13).Value = 52.441971
$ws.Cells.Item(11, 14).Value = 157.325913
$ws.Cells.Item(11, 15).Value = 0.2450515065683088
$ws.Cells.Item(11, 16).Value = 0.2450515065683087
$ws.Cells.Item(11, 17).Value = 603.7825570449661
$ws.Cells.Item(11, 18).Value = 5434.043013404695
$ws.Cells.Item(11, 19).Value = 0.04404787611115974
$ws.Cells.Item(11, 20).Value = 0.04404787611115974

# Row 12
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 7).Value = 11.513346
$ws.Cells.Item(12, 8).Value = 34.540038
$ws.Cells.Item(12, 9).Value = 0.179749460544048
$ws.Cells.Item(12, 10).Value = 0.179749460544048
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 13).Value = 54.667459
$ws.Cells.Item(12, 14).Value = 164.002377
$ws.Cells.Item(12, 15).Value = 0.255450795093328
$ws.Cells.Item(12, 16).Value = 0.255450795093328
$ws.Cells.Item(12, 17).Value = 629.4053704078141
$ws.Cells.Item(12, 18).Value = 5664.648333670326
$ws.Cells.Item(12, 19).Value = 0.04591714261357386
$ws.Cells.Item(12, 20).Value = 0.04591714261357385

# Row 13
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 7).Value = 11.513346
$ws.Cells.Item(13, 8).Value = 34.540038
$ws.Cells.Item(13, 9).Value = 0.179749460544048
$ws.Cells.Item(13, 10).Value = 0.179749460544048
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 13).Value = 18.95316166666667
$ws.Cells.Item(13, 14).Value = 56.85948500000001
$ws.Cells.Item(13, 15).Value = 0.08856457398691947
$ws.Cells.Item(13, 16).Value = 0.08856457398691944
$ws.Cells.Item(13, 17).Value = 218.21430806227
$ws.Cells.Item(13, 18).Value = 1963.92877256043
$ws.Cells.Item(13, 19).Value = 0.01591943439746221
$ws.Cells.Item(13, 20).Value = 0.0159194343974622

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 7).Value = 5.932994999999999
$ws.Cells.Item(14, 8).Value = 17.798985
$ws.Cells.Item(14, 9).Value = 0.09262751685396531
$ws.Cells.Item(14, 10).Value = 0.09262751685396531
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 13).Value = 87.94127933333334
$ws.Cells.Item(14, 14).Value = 263.823838
$ws.Cells.Item(14, 15).Value = 0.4109331243514438
$ws.Cells.Item(14, 16).Value = 0.4109331243514437
$ws.Cells.Item(14, 17).Value = 521.7551705782699
$ws.Cells.Item(14, 18).Value = 4695.79653520443
$ws.Cells.Item(14, 19).Value = 0.03806371490171599
$ws.Cells.Item(14, 20).Value = 0.03806371490171598

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 7).Value = 5.932994999999999
$ws.Cells.Item(15, 8).Value = 17.798985
$ws.Cells.Item(15, 9).Value = 0.09262751685396531
$ws.Cells.Item(15, 10).Value = 0.09262751685396531
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 13).Value = 52.441971
$ws.Cells.Item(15, 14).Value = 157.325913
$ws.Cells.Item(15, 15).Value = 0.2450515065683088
$ws.Cells.Item(15, 16).Value = 0.2450515065683087
$ws.Cells.Item(15, 17).Value = 311.137951733145
$ws.Cells.Item(15, 18).Value = 2800.241565598305
$ws.Cells.Item(15, 19).Value = 0.02269851255474561
$ws.Cells.Item(15, 20).Value = 0.02269851255474561

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 7).Value = 5.932994999999999
$ws.Cells.Item(16, 8).Value = 17.798985
$ws.Cells.Item(16, 9).Value = 0.09262751685396531
$ws.Cells.Item(16, 10).Value = 0.09262751685396531
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 13).Value = 54.667459
$ws.Cells.Item(16, 14).Value = 164.002377
$ws.Cells.Item(16, 15).Value = 0.255450795093328
$ws.Cells.Item(16, 16).Value = 0.255450795093328
$ws.Cells.Item(16, 17).Value = 324.3417609097049
$ws.Cells.Item(16, 18).Value = 2919.075848187344
$ws.Cells.Item(16, 19).Value = 0.02366177282786608
$ws.Cells.Item(16, 20).Value = 0.02366177282786608

# Row 17
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 7).Value = 5.932994999999999
$ws.Cells.Item(17, 8).Value = 17.798985
$ws.Cells.Item(17, 9).Value = 0.09262751685396531
$ws.Cells.Item(17, 10).Value = 0.09262751685396531
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 13).Value = 18.95316166666667
$ws.Cells.Item(17, 14).Value = 56.85948500000001
$ws.Cells.Item(17, 15).Value = 0.08856457398691947
$ws.Cells.Item(17, 16).Value = 0.08856457398691944
$ws.Cells.Item(17, 17).Value = 112.449013402525
$ws.Cells.Item(17, 18).Value = 1012.041120622725
$ws.Cells.Item(17, 19).Value = 0.008203516569637638
$ws.Cells.Item(17, 20).Value = 0.008203516569637638
